# feat: add 2022-Q4 data
#
# 1) Insert a new "2022-Q4" worksheet (fund-holding detail), positioned
#    right after "总计" and before "2022-Q3" - built as a copy of the
#    "2022-Q3" sheet (so header row + styles match exactly), then its
#    two data rows are overwritten with the new quarter's fund data.
# 2) Update the "总计" (total) summary sheet: insert a new top data row
#    for 2022-Q4, pushing the existing quarterly rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: new "2022-Q4" fund-detail sheet
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")

# Worksheet.Copy(Before) duplicates data + styles + column widths and
# places the duplicate immediately before $wsQ3 - exactly where the new
# quarter tab belongs (right after "总计").
$wsQ3.Copy($wsQ3)
$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"

# The fund-code / metric columns hold numeric-looking text (e.g. a
# leading-zero code like "009439" or a percentage string) in the source
# data, so force text formatting before assigning - otherwise Excel's
# smart entry would coerce them to numbers and drop the leading zero.
$wsQ4.Range("B2:B3").NumberFormat = "@"
$wsQ4.Range("D2:G3").NumberFormat = "@"

$wsQ4.Cells.Item(2, 2).Value = "009439"
$wsQ4.Cells.Item(2, 3).Value = "西部利得中证国有企业红利指数增强（LOF）C"
$wsQ4.Cells.Item(2, 4).Value = "5.98"
$wsQ4.Cells.Item(2, 5).Value = "89.71"
$wsQ4.Cells.Item(2, 6).Value = "2.36"
$wsQ4.Cells.Item(2, 7).Value = "0.1411"
$wsQ4.Cells.Item(2, 8).Value = 9

$wsQ4.Cells.Item(3, 2).Value = "501059"
$wsQ4.Cells.Item(3, 3).Value = "西部利得中证国有企业红利指数增强（LOF）A"
$wsQ4.Cells.Item(3, 4).Value = "5.18"
$wsQ4.Cells.Item(3, 5).Value = "89.71"
$wsQ4.Cells.Item(3, 6).Value = "2.36"
$wsQ4.Cells.Item(3, 7).Value = "0.1222"
$wsQ4.Cells.Item(3, 8).Value = 9

# ---------------------------------------------------------------------
# Part 2: "总计" summary sheet - add the 2022-Q4 row on top, shifting
# the previously existing rows down by one (their own values are
# unchanged, only their row position / running index moves).
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$quarters = @(
    @("2022-Q4", 2, 0.26),
    @("2022-Q3", 2, 0.06),
    @("2022-Q2", 1, 0.03),
    @("2022-Q1", 4, 0.29),
    @("2021-Q4", 2, 0.48),
    @("2021-Q3", 2, 0.03),
    @("2021-Q2", 1, 0.02)
)

# Row 8 is brand new (beyond the original A1:D7 extent) - give its "A"
# cell the same style as the existing index cells above (bold/centered,
# style index shared by A2:A7) before filling in the values.
$wsTotal.Range("A7").Copy()
$wsTotal.Range("A8").PasteSpecial(-4122)

for ($i = 0; $i -lt $quarters.Count; $i++) {
    $row = $i + 2
    $wsTotal.Cells.Item($row, 1).Value = $i
    $wsTotal.Cells.Item($row, 2).Value = $quarters[$i][0]
    $wsTotal.Cells.Item($row, 3).Value = $quarters[$i][1]
    $wsTotal.Cells.Item($row, 4).Value = $quarters[$i][2]
}
